# fix: eager set cast of possibly None implementation_groups field
#
# Adds a "library_publication_date" row to the library_content sheet
# (sheet1), bumps the library_version value, and marks the "mesures"
# sheet row for POL.ACCESS requirement (row 44) with implementation
# group "S" in column C — fixing the missing IG that caused a null set
# when implementation_groups was eagerly cast.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # library_content
$ws2 = $wb.Worksheets.Item(2)   # mesures

# --- library_content (sheet1) ---------------------------------------

# library_version bumped from 2 to 3
$ws1.Range("B2").Value = 3

# Insert a new row 4 for "library_publication_date" (pushes the
# existing rows 4-20 down to 5-21).
$ws1.Rows("4:4").Insert()
$ws1.Range("A4").Value = "library_publication_date"
$ws1.Range("B4").Value = 45680
$ws1.Range("B4").NumberFormat = "mm-dd-yy"
$ws1.Range("B4").HorizontalAlignment = -4131

# --- mesures (sheet2) -------------------------------------------------

# Row 44 (POL.ACCESS requirement) was missing its implementation group;
# set it to "S" in column C, matching the style of the rest of the row.
$ws2.Range("C44").Value = "S"
$ws2.Range("C44").VerticalAlignment = -4108
$ws2.Range("C44").WrapText = $true

# --- view/selection state ---------------------------------------------

# Scroll mesures to the area of interest and leave library_content as
# the active/selected tab when the workbook is saved.
$ws2.Activate()
$ws2.Range("E12").Select()

$ws1.Activate()
$ws1.Range("A4").Select()
